$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "27.505.67"
$ws.Range("E2").Value = "  -0.70%  "
$ws.Range("D3").Value = "1.620.75"
$ws.Range("E3").Value = "  -1.48%  "
$ws.Range("E4").Value = "  +0.06%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.522"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -0.83%  "
$ws.Range("E7").Value = "  +0.06%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "23.02"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -1.62%  "
$ws.Range("E9").Value = "  +1.58%  "
$ws.Range("E10").Value = "  -0.18%  "
$ws.Range("E11").Value = "  -1.71%  "
$ws.Range("D12").Value = "1.851.47"
$ws.Range("E12").Value = "  -1.40%  "
$ws.Range("D13").Value = "1.647.42"
$ws.Range("E13").Value = "  +0.11%  "
$ws.Range("E15").Value = "  -2.29%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "65.16"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +0.61%  "
$ws.Range("D17").Value = "27.488.01"
$ws.Range("E17").Value = "  -0.72%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "229.67"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -0.72%  "
$ws.Range("D19").Value = "0.0₃0717"
$ws.Range("E19").Value = "  -1.04%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "7.54"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -2.06%  "
$ws.Range("E21").Value = "  +0.05%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "10.38"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +3.23%  "
$ws.Range("E23").Value = "  +0.53%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.12"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +7.98%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "149.24"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -0.48%  "
$ws.Range("E26").Value = "  -1.30%  "
$ws.Range("E27").Value = "  -0.82%  "
$ws.Range("E28").Value = "  +0.00%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "15.52"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -1.02%  "
$ws.Range("E30").Value = "  -0.88%  "
$ws.Range("E31").Value = "  -0.99%  "
$ws.Range("E32").Value = "  -1.52%  "
$ws.Range("D33").Value = "1.463.58"
$ws.Range("E33").Value = "  +0.79%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "3.05"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -2.82%  "
$ws.Range("E35").Value = "  -2.68%  "
$ws.Range("E36").Value = "  -0.43%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.944"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +6.50%  "
$ws.Range("E38").Value = "  -0.22%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.870"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -1.64%  "
$ws.Range("E40").Value = "  -3.23%  "
$ws.Range("E41").Value = "  +0.04%  "
$ws.Range("E42").Value = "  -2.69%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "67.24"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -5.48%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "2.45"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -0.78%  "
$ws.Range("E45").Value = "  -2.06%  "
$ws.Range("E46").Value = "  -6.30%  "
$ws.Range("E47").Value = "  +1.78%  "
$ws.Range("D48").Value = "1.760.74"
$ws.Range("E48").Value = "  -1.53%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "87.10"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +1.31%  "
$ws.Range("E50").Value = "  -1.23%  "
